$d = $word.ActiveDocument

# --- 1. Body: "TERE" (bold run, after "A ") -> "QWER" ---
$d.Content.Find.Execute("TERE", $true, $false, $false, $false, $false, `
                         $true, 0, $false, "QWER", 1) | Out-Null

# --- 2. Header: several case-sensitive replacements, applied in document order ---
$hdr = $d.Sections.Item(1).Headers.Item(1)

$script:hdrPos = 0

function Replace-NextMatch($pattern, $replacement) {
    $rng = $hdr.Range.Duplicate
    $rng.Start = $script:hdrPos
    $rng.End = $hdr.Range.End
    $found = $rng.Find.Execute($pattern, $true, $false, $false, $false, $false, `
                                $true, 0, $false, $replacement, 1)
    if ($found) {
        $script:hdrPos = $rng.End
    }
    return $found
}

# "DIRETORIA DE ENSINO REGIAO TRE" -> "...QWER"
Replace-NextMatch "TRE" "QWER" | Out-Null

# standalone "TERE" before " - DEP." -> "QWER"
Replace-NextMatch "TERE" "QWER" | Out-Null

# five "Tre" runs, in document order
Replace-NextMatch "Tre" "Qwer" | Out-Null
Replace-NextMatch "Tre" "Qwer" | Out-Null
Replace-NextMatch "Tre" "Qewr" | Out-Null
Replace-NextMatch "Tre" "Qewr" | Out-Null
Replace-NextMatch "Tre" "Qwer" | Out-Null

# three "tre" runs (CEP / Tel / Email), in document order
Replace-NextMatch "tre" "qwer" | Out-Null
Replace-NextMatch "tre" "qwer" | Out-Null
Replace-NextMatch "tre" "qwer" | Out-Null
